{"js": "// Replace each arithmetic-exercise answer in the table with its updated\n// value. The mapping below is old-text -> new-text, built from the\n// authoritative OOXML diff; every old value is unique in the document so a\n// straight text lookup (applied once per paragraph, in document order) is\n// unambiguous and reproduces the diff exactly.\nconst replacements = [[\"57+35=92\",\"59+7=66\"],[\"59+12=71\",\"47+7=54\"],[\"90-66=24\",\"19+64=83\"],[\"55+8=63\",\"77+4=81\"],[\"93-85=8\",\"75-36=39\"],[\"75+17=92\",\"17+79=96\"],[\"58+4=62\",\"57+36=93\"],[\"90-71=19\",\"48-9=39\"],[\"16+79=95\",\"5+67=72\"],[\"78+16=94\",\"59+3=62\"],[\"54-45=9\",\"81-75=6\"],[\"41-27=14\",\"57+26=83\"],[\"74-15=59\",\"78+13=91\"],[\"70-23=47\",\"20-15=5\"],[\"95-26=69\",\"2+59=61\"],[\"62+29=91\",\"66-59=7\"],[\"10-5=5\",\"65-27=38\"],[\"14+79=93\",\"83-78=5\"],[\"97-69=28\",\"53-24=29\"],[\"81-7=74\",\"81-18=63\"],[\"93-74=19\",\"6+86=92\"],[\"7+19=26\",\"63-25=38\"],[\"19+38=57\",\"50-9=41\"],[\"44+29=73\",\"14+78=92\"],[\"38+15=53\",\"23+59=82\"],[\"16+68=84\",\"42-13=29\"],[\"50-41=9\",\"28+45=73\"],[\"14+7=21\",\"18+28=46\"],[\"6+56=62\",\"77+17=94\"],[\"37+26=63\",\"73+18=91\"],[\"82-39=43\",\"64+8=72\"],[\"60-32=28\",\"40-29=11\"],[\"92-3=89\",\"62-58=4\"],[\"39+15=54\",\"6+78=84\"],[\"7+26=33\",\"54-38=16\"],[\"29+47=76\",\"50-19=31\"],[\"97-19=78\",\"93-88=5\"],[\"60-42=18\",\"97-39=58\"],[\"82-64=18\",\"33+58=91\"],[\"4+38=42\",\"98-89=9\"],[\"47+6=53\",\"76+9=85\"],[\"90-25=65\",\"17+74=91\"],[\"53-15=38\",\"44-6=38\"],[\"83-56=27\",\"71-28=43\"],[\"77-49=28\",\"98-39=59\"],[\"18+9=27\",\"41-8=33\"],[\"79+9=88\",\"61-38=23\"],[\"92-55=37\",\"80-28=52\"],[\"86-39=47\",\"47+24=71\"],[\"38+59=97\",\"7+79=86\"],[\"18+27=45\",\"50-27=23\"],[\"81-76=5\",\"72-63=9\"],[\"34+17=51\",\"77+19=96\"],[\"46+19=65\",\"59+35=94\"],[\"27-9=18\",\"42-8=34\"],[\"92-33=59\",\"9+23=32\"],[\"83-58=25\",\"91-55=36\"],[\"8+39=47\",\"74-38=36\"],[\"61-14=47\",\"66-58=8\"],[\"72-28=44\",\"39+45=84\"],[\"40-12=28\",\"97-38=59\"],[\"29+6=35\",\"47+19=66\"],[\"54-46=8\",\"8+37=45\"],[\"97-88=9\",\"40-22=18\"],[\"86-7=79\",\"27-18=9\"],[\"62-9=53\",\"8+78=86\"],[\"7+38=45\",\"31-17=14\"],[\"25+49=74\",\"25+39=64\"],[\"80-76=4\",\"64+8=72\"],[\"71-42=29\",\"9+83=92\"],[\"86-17=69\",\"64-39=25\"],[\"40-36=4\",\"18+47=65\"],[\"70-51=19\",\"59+32=91\"],[\"86-28=58\",\"28+37=65\"],[\"90-77=13\",\"45-38=7\"],[\"88-39=49\",\"63+28=91\"],[\"14+18=32\",\"19+49=68\"],[\"75+9=84\",\"17+74=91\"],[\"34+57=91\",\"8+55=63\"],[\"5+27=32\",\"22-18=4\"],[\"19+12=31\",\"76+16=92\"],[\"86-27=59\",\"8+8=16\"],[\"58+8=66\",\"80-39=41\"],[\"90-79=11\",\"90-47=43\"],[\"42+19=61\",\"57+17=74\"],[\"91-33=58\",\"50-49=1\"],[\"25+7=32\",\"51-18=33\"],[\"81-55=26\",\"7+77=84\"],[\"48+45=93\",\"9+53=62\"],[\"61-53=8\",\"35+37=72\"],[\"39+46=85\",\"72-24=48\"],[\"49+5=54\",\"73-17=56\"],[\"58+24=82\",\"18+79=97\"],[\"26+55=81\",\"63+8=71\"],[\"12-4=8\",\"71-48=23\"],[\"84-8=76\",\"32+9=41\"],[\"50-34=16\",\"64-28=36\"],[\"84-27=57\",\"19+52=71\"],[\"56+25=81\",\"77+5=82\"],[\"35-8=27\",\"19+53=72\"]];\n\n// Build a lookup map; every key is unique in this document.\nconst map = new Map(replacements);\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items.forEach(p => p.load(\"text\"));\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  const current = p.text;\n  if (map.has(current)) {\n    p.insertText(map.get(current), Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each arithmetic-exercise answer in the table with its updated\n# value. The mapping below is old-text -> new-text, built from the\n# authoritative OOXML diff; every old value is unique in the document (and\n# none collides with any replacement's new value), so running Find/Replace\n# All once per pair, in any order, reproduces the diff exactly.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old=\"57+35=92\"; New=\"59+7=66\"},\n    @{Old=\"59+12=71\"; New=\"47+7=54\"},\n    @{Old=\"90-66=24\"; New=\"19+64=83\"},\n    @{Old=\"55+8=63\"; New=\"77+4=81\"},\n    @{Old=\"93-85=8\"; New=\"75-36=39\"},\n    @{Old=\"75+17=92\"; New=\"17+79=96\"},\n    @{Old=\"58+4=62\"; New=\"57+36=93\"},\n    @{Old=\"90-71=19\"; New=\"48-9=39\"},\n    @{Old=\"16+79=95\"; New=\"5+67=72\"},\n    @{Old=\"78+16=94\"; New=\"59+3=62\"},\n    @{Old=\"54-45=9\"; New=\"81-75=6\"},\n    @{Old=\"41-27=14\"; New=\"57+26=83\"},\n    @{Old=\"74-15=59\"; New=\"78+13=91\"},\n    @{Old=\"70-23=47\"; New=\"20-15=5\"},\n    @{Old=\"95-26=69\"; New=\"2+59=61\"},\n    @{Old=\"62+29=91\"; New=\"66-59=7\"},\n    @{Old=\"10-5=5\"; New=\"65-27=38\"},\n    @{Old=\"14+79=93\"; New=\"83-78=5\"},\n    @{Old=\"97-69=28\"; New=\"53-24=29\"},\n    @{Old=\"81-7=74\"; New=\"81-18=63\"},\n    @{Old=\"93-74=19\"; New=\"6+86=92\"},\n    @{Old=\"7+19=26\"; New=\"63-25=38\"},\n    @{Old=\"19+38=57\"; New=\"50-9=41\"},\n    @{Old=\"44+29=73\"; New=\"14+78=92\"},\n    @{Old=\"38+15=53\"; New=\"23+59=82\"},\n    @{Old=\"16+68=84\"; New=\"42-13=29\"},\n    @{Old=\"50-41=9\"; New=\"28+45=73\"},\n    @{Old=\"14+7=21\"; New=\"18+28=46\"},\n    @{Old=\"6+56=62\"; New=\"77+17=94\"},\n    @{Old=\"37+26=63\"; New=\"73+18=91\"},\n    @{Old=\"82-39=43\"; New=\"64+8=72\"},\n    @{Old=\"60-32=28\"; New=\"40-29=11\"},\n    @{Old=\"92-3=89\"; New=\"62-58=4\"},\n    @{Old=\"39+15=54\"; New=\"6+78=84\"},\n    @{Old=\"7+26=33\"; New=\"54-38=16\"},\n    @{Old=\"29+47=76\"; New=\"50-19=31\"},\n    @{Old=\"97-19=78\"; New=\"93-88=5\"},\n    @{Old=\"60-42=18\"; New=\"97-39=58\"},\n    @{Old=\"82-64=18\"; New=\"33+58=91\"},\n    @{Old=\"4+38=42\"; New=\"98-89=9\"},\n    @{Old=\"47+6=53\"; New=\"76+9=85\"},\n    @{Old=\"90-25=65\"; New=\"17+74=91\"},\n    @{Old=\"53-15=38\"; New=\"44-6=38\"},\n    @{Old=\"83-56=27\"; New=\"71-28=43\"},\n    @{Old=\"77-49=28\"; New=\"98-39=59\"},\n    @{Old=\"18+9=27\"; New=\"41-8=33\"},\n    @{Old=\"79+9=88\"; New=\"61-38=23\"},\n    @{Old=\"92-55=37\"; New=\"80-28=52\"},\n    @{Old=\"86-39=47\"; New=\"47+24=71\"},\n    @{Old=\"38+59=97\"; New=\"7+79=86\"},\n    @{Old=\"18+27=45\"; New=\"50-27=23\"},\n    @{Old=\"81-76=5\"; New=\"72-63=9\"},\n    @{Old=\"34+17=51\"; New=\"77+19=96\"},\n    @{Old=\"46+19=65\"; New=\"59+35=94\"},\n    @{Old=\"27-9=18\"; New=\"42-8=34\"},\n    @{Old=\"92-33=59\"; New=\"9+23=32\"},\n    @{Old=\"83-58=25\"; New=\"91-55=36\"},\n    @{Old=\"8+39=47\"; New=\"74-38=36\"},\n    @{Old=\"61-14=47\"; New=\"66-58=8\"},\n    @{Old=\"72-28=44\"; New=\"39+45=84\"},\n    @{Old=\"40-12=28\"; New=\"97-38=59\"},\n    @{Old=\"29+6=35\"; New=\"47+19=66\"},\n    @{Old=\"54-46=8\"; New=\"8+37=45\"},\n    @{Old=\"97-88=9\"; New=\"40-22=18\"},\n    @{Old=\"86-7=79\"; New=\"27-18=9\"},\n    @{Old=\"62-9=53\"; New=\"8+78=86\"},\n    @{Old=\"7+38=45\"; New=\"31-17=14\"},\n    @{Old=\"25+49=74\"; New=\"25+39=64\"},\n    @{Old=\"80-76=4\"; New=\"64+8=72\"},\n    @{Old=\"71-42=29\"; New=\"9+83=92\"},\n    @{Old=\"86-17=69\"; New=\"64-39=25\"},\n    @{Old=\"40-36=4\"; New=\"18+47=65\"},\n    @{Old=\"70-51=19\"; New=\"59+32=91\"},\n    @{Old=\"86-28=58\"; New=\"28+37=65\"},\n    @{Old=\"90-77=13\"; New=\"45-38=7\"},\n    @{Old=\"88-39=49\"; New=\"63+28=91\"},\n    @{Old=\"14+18=32\"; New=\"19+49=68\"},\n    @{Old=\"75+9=84\"; New=\"17+74=91\"},\n    @{Old=\"34+57=91\"; New=\"8+55=63\"},\n    @{Old=\"5+27=32\"; New=\"22-18=4\"},\n    @{Old=\"19+12=31\"; New=\"76+16=92\"},\n    @{Old=\"86-27=59\"; New=\"8+8=16\"},\n    @{Old=\"58+8=66\"; New=\"80-39=41\"},\n    @{Old=\"90-79=11\"; New=\"90-47=43\"},\n    @{Old=\"42+19=61\"; New=\"57+17=74\"},\n    @{Old=\"91-33=58\"; New=\"50-49=1\"},\n    @{Old=\"25+7=32\"; New=\"51-18=33\"},\n    @{Old=\"81-55=26\"; New=\"7+77=84\"},\n    @{Old=\"48+45=93\"; New=\"9+53=62\"},\n    @{Old=\"61-53=8\"; New=\"35+37=72\"},\n    @{Old=\"39+46=85\"; New=\"72-24=48\"},\n    @{Old=\"49+5=54\"; New=\"73-17=56\"},\n    @{Old=\"58+24=82\"; New=\"18+79=97\"},\n    @{Old=\"26+55=81\"; New=\"63+8=71\"},\n    @{Old=\"12-4=8\"; New=\"71-48=23\"},\n    @{Old=\"84-8=76\"; New=\"32+9=41\"},\n    @{Old=\"50-34=16\"; New=\"64-28=36\"},\n    @{Old=\"84-27=57\"; New=\"19+52=71\"},\n    @{Old=\"56+25=81\"; New=\"77+5=82\"},\n    @{Old=\"35-8=27\"; New=\"19+53=72\"}\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute(\n        $pair.Old,   # FindText\n        $false,      # MatchCase\n        $false,      # MatchWholeWord\n        $false,      # MatchWildcards\n        $false,      # MatchSoundsLike\n        $false,      # MatchAllWordForms\n        $true,       # Forward\n        1,           # Wrap (wdFindContinue)\n        $false,      # Format\n        $pair.New,   # ReplaceWith\n        2            # Replace (wdReplaceAll)\n    )\n}\n"}
